# Append a new job-listing row ("貸別荘収支表自動集計システム構築の依頼") as the new
# row 6 of the "ランサーズ" sheet, push the previous row 6
# ("プログラム修正依頼!...") down to row 7, and refresh every
# "取得日時" timestamp in A2:A7 to 2025-11-24 01:25:09.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2025-11-24 01:25:09"

# This engine's Range.Hyperlinks.Delete() removes every hyperlink on the
# sheet (not just the ones in the target range), so we drop them all here
# and rebuild the full set below in the same left-to-right order. That
# recreates the relationship ids (rId1..rId6) exactly as Excel would after
# inserting a row that carries a hyperlink along with it.
$ws.Range("F2").Hyperlinks.Delete()

# Refresh the capture timestamp on the four untouched rows.
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp

# Move the former row 6 ("プログラム修正依頼") down to row 7.
$ws.Range("A7").Value = $newTimestamp
$ws.Range("B7").Value = "【急募】プログラム修正依頼!スキルを活かしてみませんか?"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5440002"
$ws.Range("G7").Value = 13

# Write the newly scraped listing into row 6.
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "【急募】貸別荘収支表自動集計システム構築の依頼"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5440042"
$ws.Range("G6").Value = 25

# Rebuild the hyperlinks: F2..F5 unchanged targets, F6 gets the brand new
# listing's URL, F7 gets the URL that used to live on (the old) F6.
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5440052")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5439921")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5440010")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5439670")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5440042")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5440002")

# Hyperlinks.Add() leaves a throwaway style behind; re-apply the workbook's
# built-in "Hyperlink" cell style everywhere so F2:F7 all share the same
# style slot they originally used.
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("F3").Style = "Hyperlink"
$ws.Range("F4").Style = "Hyperlink"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("F6").Style = "Hyperlink"
$ws.Range("F7").Style = "Hyperlink"
